$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 67: next quarterly date (2025-04-01) with same pattern/style as prior rows
$ws.Range("A67").Value = 45748
$ws.Range("A67").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("B67").Value = 0
$ws.Range("C67").Value = 0.5
$ws.Range("D67").Value = 2.5
$ws.Range("E67").Value = 3.5
$ws.Range("F67").Value = 4.5
$ws.Range("G67").Value = 9.5
$ws.Range("H67").Value = 14.5
